$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) / 1h volume change (E) figures for the crypto snapshot.
$updates = @(
    @{ Row = 2; D = "29.332.35"; E = "  -1.16%  " }
    @{ Row = 3; D = "1.894.34"; E = "  -1.40%  " }
    @{ Row = 4; D = "1.001"; E = "  +0.04%  " }
    @{ Row = 5; D = "322.95"; E = "  -3.57%  " }
    @{ Row = 6; D = "1.000"; E = "  -0.01%  " }
    @{ Row = 7; D = $null; E = "  +2.25%  " }
    @{ Row = 8; D = "0.4048"; E = "  -2.22%  " }
    @{ Row = 9; D = "0.08022"; E = "  -0.53%  " }
    @{ Row = 10; D = "0.9998"; E = "  -2.02%  " }
    @{ Row = 11; D = "23.35"; E = "  +4.73%  " }
    @{ Row = 12; D = "1.859.65"; E = "  -2.67%  " }
    @{ Row = 13; D = "5.922"; E = "  -1.64%  " }
    @{ Row = 14; D = "7.040"; E = "  -2.06%  " }
    @{ Row = 15; D = "89.36"; E = "  -0.55%  " }
    @{ Row = 16; D = "1.002"; E = "  +0.09%  " }
    @{ Row = 17; D = "0.06680"; E = "  +1.30%  " }
    @{ Row = 18; D = $null; E = "  -0.93%  " }
    @{ Row = 19; D = "17.56"; E = "  -1.55%  " }
    @{ Row = 20; D = "0.9990"; E = "  -0.10%  " }
    @{ Row = 21; D = "29.355.96"; E = "  -1.06%  " }
    @{ Row = 22; D = "5.517"; E = "  -0.72%  " }
    @{ Row = 23; D = "11.68"; E = "  +0.73%  " }
    @{ Row = 24; D = "2.154"; E = "  -2.09%  " }
    @{ Row = 25; D = "2.127.07"; E = "  -0.62%  " }
    @{ Row = 26; D = "154.03"; E = "  -2.45%  " }
    @{ Row = 27; D = "19.75"; E = "  -0.96%  " }
    @{ Row = 28; D = "5.938"; E = "  +3.51%  " }
    @{ Row = 29; D = "2.088"; E = "  -3.03%  " }
    @{ Row = 30; D = "117.77"; E = "  -0.17%  " }
    @{ Row = 31; D = "1.020"; E = "  -2.62%  " }
    @{ Row = 32; D = "0.09449"; E = "  -0.04%  " }
    @{ Row = 33; D = $null; E = "  +0.00%  " }
    @{ Row = 34; D = "1.379"; E = "  -3.81%  " }
    @{ Row = 35; D = "5.350"; E = "  -1.68%  " }
    @{ Row = 36; D = "0.02244"; E = "  -1.07%  " }
    @{ Row = 37; D = "0.06036"; E = "  -1.90%  " }
    @{ Row = 38; D = $null; E = "  -0.69%  " }
    @{ Row = 39; D = "0.5848"; E = "  -1.15%  " }
    @{ Row = 40; D = "7.837"; E = "  -7.46%  " }
    @{ Row = 41; D = "0.1835"; E = "  -0.59%  " }
    @{ Row = 42; D = "10.08"; E = "  -1.76%  " }
    @{ Row = 43; D = $null; E = "  +2.59%  " }
    @{ Row = 44; D = "2.383"; E = "  +1.87%  " }
    @{ Row = 45; D = "0.07699"; E = "  +2.34%  " }
    @{ Row = 46; D = "12.23"; E = "  -0.38%  " }
    @{ Row = 47; D = "0.5488"; E = "  -1.80%  " }
    @{ Row = 48; D = "1.918"; E = "  -1.27%  " }
    @{ Row = 49; D = "113.09"; E = "  +0.06%  " }
    @{ Row = 50; D = "0.2957"; E = "  -1.24%  " }
    @{ Row = 51; D = "43.64"; E = "  -0.73%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $dCell = $ws.Range("D$($u.Row)")
        $dCell.NumberFormat = "@"   # keep literal text, e.g. "1.000" / "29.332.35"
        $dCell.Value = $u.D
    }
    if ($u.E -ne $null) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
